$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{Row=2; B=$null; C=$null; D="23.618.83"; E="  -1.42%  "}
  @{Row=3; B=$null; C=$null; D="1.619.71"; E="  -2.24%  "}
  @{Row=4; B=$null; C=$null; D="1.004"; E="  +0.06%  "}
  @{Row=5; B=$null; C=$null; D="1.006"; E="  +0.50%  "}
  @{Row=6; B=$null; C=$null; D="306.82"; E="  -0.89%  "}
  @{Row=7; B=$null; C=$null; D="0.3833"; E="  -1.67%  "}
  @{Row=8; B=$null; C=$null; D="0.3762"; E="  -2.77%  "}
  @{Row=9; B=$null; C=$null; D="49.43"; E="  -3.61%  "}
  @{Row=10; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.010"; E="  +0.64%  "}
  @{Row=11; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="1.302"; E="  -4.40%  "}
  @{Row=12; B=$null; C=$null; D="0.08239"; E="  -2.93%  "}
  @{Row=13; B=$null; C=$null; D="23.29"; E="  -2.44%  "}
  @{Row=14; B=$null; C=$null; D="6.809"; E="  -5.46%  "}
  @{Row=15; B=$null; C=$null; D="7.641"; E="  -4.62%  "}
  @{Row=16; B=$null; C=$null; D="0.00001283"; E="  -2.12%  "}
  @{Row=17; B=$null; C=$null; D="1.620.39"; E="  -2.72%  "}
  @{Row=18; B=$null; C=$null; D="92.86"; E="  -1.78%  "}
  @{Row=19; B=$null; C=$null; D="0.06918"; E="  -0.92%  "}
  @{Row=20; B=$null; C=$null; D="19.01"; E="  -4.42%  "}
  @{Row=21; B=$null; C=$null; D=$null; E="  -2.92%  "}
  @{Row=22; B=$null; C=$null; D="1.005"; E="  +0.38%  "}
  @{Row=23; B=$null; C=$null; D="13.36"; E="  -2.06%  "}
  @{Row=24; B=$null; C=$null; D="23.621.85"; E="  -1.47%  "}
  @{Row=25; B=$null; C=$null; D="2.405"; E="  -3.35%  "}
  @{Row=26; B=$null; C=$null; D="2.863"; E="  -6.95%  "}
  @{Row=27; B=$null; C=$null; D="21.62"; E="  -2.84%  "}
  @{Row=28; B=$null; C=$null; D="152.13"; E="  -0.81%  "}
  @{Row=29; B=$null; C=$null; D="5.421"; E="  +2.29%  "}
  @{Row=30; B=$null; C=$null; D="7.911"; E="  -0.54%  "}
  @{Row=31; B=$null; C=$null; D="134.43"; E="  -3.97%  "}
  @{Row=32; B=$null; C=$null; D="2.478"; E="  -0.32%  "}
  @{Row=33; B=$null; C=$null; D="1.808.63"; E="  -2.12%  "}
  @{Row=34; B=$null; C=$null; D="0.9672"; E="  -7.14%  "}
  @{Row=35; B=$null; C=$null; D="0.07735"; E="  -4.90%  "}
  @{Row=36; B=$null; C=$null; D="0.02843"; E="  -5.02%  "}
  @{Row=37; B=$null; C=$null; D="6.485"; E="  -3.33%  "}
  @{Row=38; B=$null; C=$null; D="0.2609"; E="  -3.29%  "}
  @{Row=39; B=$null; C=$null; D="10.29"; E="  -7.03%  "}
  @{Row=40; B=$null; C=$null; D="0.08979"; E="  -1.83%  "}
  @{Row=41; B=$null; C=$null; D="0.7349"; E="  -2.58%  "}
  @{Row=42; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="13.13"; E="  -3.75%  "}
  @{Row=43; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.398"; E="  -1.48%  "}
  @{Row=44; B=$null; C=$null; D="16.11"; E="  -2.50%  "}
  @{Row=45; B=$null; C=$null; D="0.6764"; E="  -3.48%  "}
  @{Row=46; B=$null; C=$null; D="2.365"; E="  -4.79%  "}
  @{Row=47; B=$null; C=$null; D="4.045"; E="  -1.08%  "}
  @{Row=48; B=$null; C=$null; D=$null; E="  +0.34%  "}
  @{Row=49; B=$null; C=$null; D="0.08114"; E="  -2.07%  "}
  @{Row=50; B=$null; C=$null; D="131.91"; E="  -2.65%  "}
  @{Row=51; B=$null; C=$null; D="1.197"; E="  -2.96%  "}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.B) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($null -ne $u.C) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).NumberFormat = "@"
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    if ($null -ne $u.E) { $ws.Cells.Item($r, 5).Value = $u.E }
}
